$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the table with new columns M:AE ---
# 1) Propagate body style (s=2, empty) from column L across the new
#    M2:AE18 block (source is a single empty-but-styled column so Excel
#    tiles it across the full destination width).
$ws.Range("L2:L18").Copy($ws.Range("M2:AE18"))

# 2) Propagate header style (s=1) from L1 across M1:AE1 the same way.
$ws.Range("L1").Copy($ws.Range("M1:AE1"))

# M1 must stay empty (style only), just like the source diff shows.
$ws.Range("M1").Value = ""

# 3) Fill in the new header captions (N1:AE1).
$ws.Range("N1").Value = "Server 2016 Classic 6.98 McAfee"
$ws.Range("O1").Value = "Server 2016 Classic 6.98 Symantec"
$ws.Range("P1").Value = "Server 2019 Classic 6.97 McAfee"
$ws.Range("Q1").Value = "Server 2019 Classic 6.97 Symantec"
$ws.Range("R1").Value = "Server 2016 Evo 24.1 McAfee"
$ws.Range("S1").Value = "Server 2016 Evo 24.1 Symantec"
$ws.Range("T1").Value = "Server 2016 Evo 24.1 Windows Defender"
$ws.Range("U1").Value = "Server 2019 Evo 23.1 McAfee"
$ws.Range("V1").Value = "Server 2019 Evo 23.1 Symantec"
$ws.Range("W1").Value = "Server 2019 Evo 23.2 McAfee"
$ws.Range("X1").Value = "Server 2019 Evo 23.2 Symantec"
$ws.Range("Y1").Value = "Server 2019 Evo 23.2 Windows Defender"
$ws.Range("Z1").Value = "Server 2022 Evo 22.2 McAfee"
$ws.Range("AA1").Value = "Server 2022 Evo 22.2 Symantec"
$ws.Range("AB1").Value = "Wind10 21H2 Classic 6.98 Evo versions: 22.X - 23.1"
$ws.Range("AC1").Value = "Win10 22H2 Classic 6.97 Evo 23.2"
$ws.Range("AD1").Value = "Win11 22H2 Evo 23.1"
$ws.Range("AE1").Value = "Win11 23H2 Evo 24.1"

# --- Fix up the "version RTM" -> "RTM" text in the Model/Version column ---
$ws.Range("D15").Value = "RTM"
$ws.Range("D16").Value = "RTM"
